# Cron / DB sync update: refresh the first contact row (row 5) on Sheet1
# with the latest data pulled from the database. The underlying shared
# strings for the name and e-mail of "Erik Ladnak" changed slightly
# (the DB export truncated a couple of characters), while everything
# else on the sheet (formatting, other rows, other sheets) stays as is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("A5").Value = "Erik Ladňá`nJakub Rončák`nElf René"
$ws.Range("B5").Value = "ladnak.erik@gmail.co`njakub.roncak@gmail.com`nelf@rene.sk"

# Leave the cursor where the editor left it after the refresh.
$ws.Range("B6").Select()
